# "further cleaning to metadata"
#  - shared string "E7760" -> "E7420" (the s2cDNASampleNumber value used in col G)
#  - column G data cells (G2:G27) get a new font (Arial 11, black) applied
#  - column H data cells (H2:H27) become formula cells "=FALSE()" (was a
#    hard-coded boolean FALSE literal)
#  - the sheet's active selection moves from H2:H27 to G2:G27

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 27

# --- update the sample-number text shown in column G ------------------
# every G2:G27 cell shares the same string ("E7760"); re-writing all of
# them to the new text keeps the shared-string table de-duplicated
# instead of leaving a stray unused entry behind.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = "E7420"
}

# --- re-style column G (sample number) with the new font --------------
# (set Size before Name/Color so the engine folds this into a single new
#  font entry instead of leaving stray intermediate fonts behind - the
#  G2:G27 cells start out on two different inherited fonts (Arial 10 /
#  Calibri 12) so order matters here)
$gFont = $ws.Range("G2:G27").Font
$gFont.Size = 11
$gFont.Name = "Arial"
$gFont.Color = 0

# --- column H becomes a FALSE() formula instead of a literal boolean --
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Formula = "=FALSE()"
}

# --- move the active selection to the (now current) data column -------
$ws.Range("G2:G27").Select() | Out-Null
